$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Required Loop" value (A13): 0.0017 -> 0.001
$ws.Range("A13").Value = 0.001

# Update the F column formulas to reference $A$5 instead of $A$6
$ws.Range("F2").Formula = "=`$A`$5*C2"
$ws.Range("F3:F9").Formula = "=`$A`$5*C3"

# Move the active selection to H8
$ws.Range("H8").Select()
